# Apply the TestData sheet updates described in the commit "Test data updated 4/21"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Unhide columns P (16) and Q (17)
$ws.Columns.Item(16).Hidden = $false
$ws.Columns.Item(17).Hidden = $false

# Swap the values of G2 and G3 (Supplier Order code)
$ws.Range("G2").Value = "BA"
$ws.Range("G3").Value = "HO"

# Fix the P/Q formulas which referenced a deleted column (#REF!) to use the
# correct source columns instead.
$ws.Range("P2").Formula = '=CONCATENATE(C2,D2,TEXT(K2,"mm/dd/yy"),O2)'
$ws.Range("Q2").Formula = '=CONCATENATE(P2,G2)'

$ws.Range("P3").Formula = '=CONCATENATE(C3,D3,TEXT(K3,"mm/dd/yy"),O3)'
$ws.Range("Q3").Formula = '=CONCATENATE(P3,G3)'

# Move the active selection to G7 as in the source workbook
$ws.Range("G7").Select()
